# This script re-applies a set of row-content rotations within the single
# worksheet of the workbook. Several groups of rows (identified by row
# number) had their entire row contents (columns A:AY) cyclically rotated
# among themselves. We:
#   1. Snapshot the full A:AY value-array of every row involved BEFORE any
#      writes happen (required because the rotations are circular).
#   2. Write each snapshot into its destination row.
#   3. Because columns Y (Startdatum) and AA (Slutdatum) hold date-like
#      text such as "2023-09-03" that Excel would otherwise auto-convert
#      into a date serial number on write, we temporarily force those two
#      columns to Text number format for the destination rows, write the
#      values, then clear the formatting again so no stray style is left
#      behind (matching the original workbook, where no cell carries an
#      explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"

# Each inner array is one rotation cycle: the data that starts in
# $cycle[0] ends up in $cycle[-1], the data in $cycle[1] ends up in
# $cycle[0], etc. In other words, for every position i in the cycle,
# the NEW content of row $cycle[i] is the OLD content of row $cycle[i+1]
# (wrapping around).
$cycles = @(
    ,@(26, 27, 28)
    ,@(53, 54)
    ,@(59, 63, 62, 61, 60)
    ,@(73, 75)
    ,@(74, 76)
)

foreach ($cycle in $cycles) {
    $n = $cycle.Length

    # Step 1: snapshot every row's full A:AY contents before writing.
    $snapshots = @{}
    foreach ($r in $cycle) {
        $rng = $ws.Range("A" + $r + ":" + $lastCol + $r)
        $snapshots[$r] = $rng.Value()
    }

    # Step 2: write rotated snapshots into destination rows.
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $n]
        $data = $snapshots[$srcRow]

        $destRange = $ws.Range("A" + $destRow + ":" + $lastCol + $destRow)
        $dateRange = $ws.Range("Y" + $destRow + ":Y" + $destRow + ",AA" + $destRow + ":AA" + $destRow)

        $ws.Range("Y" + $destRow).NumberFormat = "@"
        $ws.Range("AA" + $destRow).NumberFormat = "@"

        $destRange.Value = $data

        $ws.Range("Y" + $destRow).ClearFormats()
        $ws.Range("AA" + $destRow).ClearFormats()
    }
}
